$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update Price (D) and Volume(1h) (E) values per row ---
# Leading apostrophe forces text interpretation for numeric-looking values;
# resetting Style to Normal afterward avoids Excel auto-applying a text number format/style.
$ws.Range("D2").Value = "'69.474.36"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +0.18%  "

$ws.Range("D3").Value = "'3.671.24"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  -0.46%  "

$ws.Range("E4").Value = "  +0.07%  "

$ws.Range("D5").Value = "'646.45"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -5.28%  "

$ws.Range("D6").Value = "'159.58"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +0.00%  "

$ws.Range("E7").Value = "  -0.01%  "

$ws.Range("E8").Value = "  +0.40%  "

$ws.Range("E9").Value = "  -0.83%  "

$ws.Range("E10").Value = "  -0.40%  "

$ws.Range("E11").Value = "  -0.18%  "

$ws.Range("E12").Value = "  -0.34%  "

$ws.Range("D13").Value = "'4.290.56"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -0.41%  "

$ws.Range("D14").Value = "'32.62"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +0.36%  "

$ws.Range("D15").Value = "'3.670.42"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -0.60%  "

$ws.Range("D16").Value = "'69.442.17"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +0.17%  "

$ws.Range("E17").Value = "  +0.69%  "

$ws.Range("D18").Value = "'15.97"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -0.62%  "

$ws.Range("E19").Value = "  -0.33%  "

$ws.Range("D20").Value = "'464.78"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -0.99%  "

$ws.Range("D21").Value = "'9.77"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -1.75%  "

$ws.Range("E22").Value = "  -1.67%  "

$ws.Range("D23").Value = "'79.46"
$ws.Range("D23").Style = "Normal"

$ws.Range("D24").Value = "'3.817.96"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -0.41%  "

$ws.Range("D25").Value = "'0.999"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -0.04%  "

$ws.Range("E26").Value = "  +1.05%  "

$ws.Range("E27").Value = "  -1.49%  "

$ws.Range("E28").Value = "  -2.58%  "

$ws.Range("E29").Value = "  -2.95%  "

$ws.Range("E30").Value = "  -4.06%  "

$ws.Range("E31").Value = "  -0.31%  "

$ws.Range("D32").Value = "'0.999"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +0.00%  "

$ws.Range("E33").Value = "  -2.87%  "

$ws.Range("E34").Value = "  -1.31%  "

$ws.Range("E35").Value = "  +3.77%  "

$ws.Range("D36").Value = "'3.660.94"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -0.39%  "

$ws.Range("E37").Value = "  +0.94%  "

$ws.Range("E38").Value = "  +0.00%  "

$ws.Range("E39").Value = "  -6.10%  "

$ws.Range("D40").Value = "'178.34"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +4.75%  "

$ws.Range("E41").Value = "  +0.03%  "

$ws.Range("D44").Value = "'0.927"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -1.76%  "

$ws.Range("D45").Value = "'46.63"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -2.11%  "

$ws.Range("D46").Value = "'2.71"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +0.02%  "

$ws.Range("E47").Value = "  -3.49%  "

$ws.Range("D48").Value = "'26.88"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -5.50%  "

$ws.Range("E49").Value = "  -3.87%  "

$ws.Range("E50").Value = "  +0.17%  "

$ws.Range("E51").Value = "  -6.54%  "

# --- Swap rows 42 and 43 (Stacks <-> Hedera) ---
$ws.Range("B42").Value = "Stacks"
$ws.Range("C42").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D42").Value = "'2.18"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -3.92%  "

$ws.Range("B43").Value = "Hedera"
$ws.Range("C43").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D43").Value = "'0.0894"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -1.49%  "
